# Product Backlog (#51): add new backlog story "Ressourcen Übersicht"
# (overall-calendar) as row 9, then restore the UI selection state that
# Excel recorded when the workbook was saved.

$wb = $excel.ActiveWorkbook

# --- Product Backlog sheet: append the new story as row 9 -----------------
$productBacklog = $wb.Sheets.Item("Product Backlog")

$productBacklog.Cells.Item(9, 1).Value = 8
$productBacklog.Cells.Item(9, 2).Value = "Ressourcen Übersicht"
$productBacklog.Cells.Item(9, 3).Value = "Auf einem Kalender sollen alle MA Einsätze dargestellt werden"
$productBacklog.Cells.Item(9, 4).Value = "low"
$productBacklog.Cells.Item(9, 5).Value = 25
$productBacklog.Cells.Item(9, 6).Value = 0
$productBacklog.Cells.Item(9, 7).Value = 0
$productBacklog.Cells.Item(9, 8).Value = "waiting"
# Match the "waiting"/"done"/... status column formatting used by the other rows
$productBacklog.Cells.Item(9, 8).VerticalAlignment = -4160

# Recorded selection on the Product Backlog sheet after the edit
[void]$productBacklog.Range("C13").Select()

# --- Sprint Backlog sheet: it stays the active tab, selection moved -------
$sprintBacklog = $wb.Sheets.Item("Sprint Backlog")
[void]$sprintBacklog.Activate()
[void]$sprintBacklog.Range("C19").Select()
